{"js": "// Ordered list of [oldText, newText] pairs describing how each non-empty\n// table-cell paragraph's text must change, in document order. The same\n// text occasionally repeats (e.g. \"23\u00f76=3, 5\" is both a target value for one\n// cell and the original value of a different, later cell), so replacement\n// MUST be done positionally (by walking cells in document order) rather\n// than via a global text search, otherwise a later replacement could\n// clobber a value just written by an earlier one.\nconst replacements = [\n  [\"37\u00f72=18, 1\", \"13\u00f79=1, 4\"],\n  [\"97\u00f75=19, 2\", \"55\u00f78=6, 7\"],\n  [\"35\u00f77=5, 0\", \"70\u00f77=10, 0\"],\n  [\"87\u00f73=29, 0\", \"90\u00f73=30, 0\"],\n  [\"63\u00f74=15, 3\", \"22\u00f73=7, 1\"],\n  [\"90\u00f79=10, 0\", \"22\u00f73=7, 1\"],\n  [\"19\u00f72=9, 1\", \"47\u00f79=5, 2\"],\n  [\"69\u00f72=34, 1\", \"86\u00f75=17, 1\"],\n  [\"24\u00f77=3, 3\", \"23\u00f76=3, 5\"],\n  [\"18\u00f76=3, 0\", \"74\u00f72=37, 0\"],\n  [\"93\u00f78=11, 5\", \"33\u00f74=8, 1\"],\n  [\"23\u00f76=3, 5\", \"38\u00f77=5, 3\"],\n  [\"33\u00f74=8, 1\", \"52\u00f72=26, 0\"],\n  [\"20\u00f79=2, 2\", \"30\u00f77=4, 2\"],\n  [\"34\u00f75=6, 4\", \"94\u00f79=10, 4\"],\n  [\"97\u00f79=10, 7\", \"10\u00f75=2, 0\"],\n  [\"47\u00f72=23, 1\", \"46\u00f73=15, 1\"],\n  [\"45\u00f78=5, 5\", \"86\u00f77=12, 2\"],\n  [\"42\u00f75=8, 2\", \"92\u00f75=18, 2\"],\n  [\"83\u00f77=11, 6\", \"66\u00f72=33, 0\"],\n  [\"96\u00f73=32, 0\", \"73\u00f78=9, 1\"],\n  [\"72\u00f75=14, 2\", \"29\u00f76=4, 5\"],\n  [\"59\u00f79=6, 5\", \"46\u00f75=9, 1\"],\n  [\"92\u00f79=10, 2\", \"72\u00f72=36, 0\"],\n  [\"63\u00f75=12, 3\", \"64\u00f79=7, 1\"],\n];\n\n// Collect every paragraph that lives inside a table cell, in document\n// (reading) order: table by table, row by row, cell by cell.\nconst tables = context.document.body.tables;\ntables.load(\"items\");\nawait context.sync();\n\nlet allCellParagraphs = [];\nfor (const table of tables.items) {\n  table.rows.load(\"items\");\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    row.cells.load(\"items\");\n  }\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    for (const cell of row.cells.items) {\n      cell.body.paragraphs.load(\"items\");\n    }\n  }\n}\nawait context.sync();\n\nfor (const table of tables.items) {\n  for (const row of table.rows.items) {\n    for (const cell of row.cells.items) {\n      for (const paragraph of cell.body.paragraphs.items) {\n        paragraph.load(\"text\");\n        allCellParagraphs.push(paragraph);\n      }\n    }\n  }\n}\nawait context.sync();\n\n// Only the paragraphs that actually contain text participate in the\n// numbered sequence described by `replacements`; the interleaved blank\n// rows are left untouched.\nconst nonEmptyParagraphs = allCellParagraphs.filter((p) => p.text.length > 0);\n\nif (nonEmptyParagraphs.length !== replacements.length) {\n  throw new Error(\n    \"Unexpected number of non-empty table cells: found \" +\n      nonEmptyParagraphs.length +\n      \", expected \" +\n      replacements.length\n  );\n}\n\nfor (let i = 0; i < replacements.length; i++) {\n  const [expectedOld, newText] = replacements[i];\n  const paragraph = nonEmptyParagraphs[i];\n  if (paragraph.text !== expectedOld) {\n    throw new Error(\n      \"Cell #\" +\n        i +\n        \" text mismatch: expected \" +\n        JSON.stringify(expectedOld) +\n        \" but found \" +\n        JSON.stringify(paragraph.text)\n    );\n  }\n  // Replacing the paragraph's text (rather than the whole cell body) keeps\n  // the existing run formatting (font, size) and paragraph formatting\n  // (justification) intact.\n  paragraph.insertText(newText, Word.InsertLocation.replace);\n}\n\nawait context.sync();\n", "ps1": "# Ordered list of (oldText, newText) pairs describing how each non-empty\n# table-cell paragraph's text must change, in document order. The same\n# text occasionally repeats (e.g. '23\u00f76=3, 5' is both the target value for\n# one cell and the original value of a different, later cell), so the\n# replacement MUST be done positionally (by walking the table's rows/cells\n# in document order) rather than via a global Find/Replace, otherwise a\n# later replacement could clobber a value just written by an earlier one.\n$replacements = @(\n    @('37\u00f72=18, 1', '13\u00f79=1, 4'),\n    @('97\u00f75=19, 2', '55\u00f78=6, 7'),\n    @('35\u00f77=5, 0', '70\u00f77=10, 0'),\n    @('87\u00f73=29, 0', '90\u00f73=30, 0'),\n    @('63\u00f74=15, 3', '22\u00f73=7, 1'),\n    @('90\u00f79=10, 0', '22\u00f73=7, 1'),\n    @('19\u00f72=9, 1', '47\u00f79=5, 2'),\n    @('69\u00f72=34, 1', '86\u00f75=17, 1'),\n    @('24\u00f77=3, 3', '23\u00f76=3, 5'),\n    @('18\u00f76=3, 0', '74\u00f72=37, 0'),\n    @('93\u00f78=11, 5', '33\u00f74=8, 1'),\n    @('23\u00f76=3, 5', '38\u00f77=5, 3'),\n    @('33\u00f74=8, 1', '52\u00f72=26, 0'),\n    @('20\u00f79=2, 2', '30\u00f77=4, 2'),\n    @('34\u00f75=6, 4', '94\u00f79=10, 4'),\n    @('97\u00f79=10, 7', '10\u00f75=2, 0'),\n    @('47\u00f72=23, 1', '46\u00f73=15, 1'),\n    @('45\u00f78=5, 5', '86\u00f77=12, 2'),\n    @('42\u00f75=8, 2', '92\u00f75=18, 2'),\n    @('83\u00f77=11, 6', '66\u00f72=33, 0'),\n    @('96\u00f73=32, 0', '73\u00f78=9, 1'),\n    @('72\u00f75=14, 2', '29\u00f76=4, 5'),\n    @('59\u00f79=6, 5', '46\u00f75=9, 1'),\n    @('92\u00f79=10, 2', '72\u00f72=36, 0'),\n    @('63\u00f75=12, 3', '64\u00f79=7, 1')\n)\n\n$d = $word.ActiveDocument\n$t = $d.Tables.Item(1)\n\n# Collect every non-empty cell in document order (row by row, left to\n# right within each row); the interleaved blank rows are skipped.\n$cells = @()\nforeach ($row in $t.Rows) {\n    foreach ($cell in $row.Cells) {\n        $raw = $cell.Range.Text\n        # Cell.Range.Text includes the trailing end-of-cell marker(s)\n        # (CR + BEL); strip those off before inspecting/comparing content.\n        $clean = $raw.TrimEnd([char]13, [char]7)\n        if ($clean.Length -gt 0) {\n            $cells += $cell\n        }\n    }\n}\n\nif ($cells.Count -ne $replacements.Count) {\n    throw \"Unexpected number of non-empty table cells: found $($cells.Count), expected $($replacements.Count)\"\n}\n\nfor ($i = 0; $i -lt $replacements.Count; $i++) {\n    $expectedOld = $replacements[$i][0]\n    $newText = $replacements[$i][1]\n    $cell = $cells[$i]\n    $current = $cell.Range.Text.TrimEnd([char]13, [char]7)\n    if ($current -ne $expectedOld) {\n        throw \"Cell #$i text mismatch: expected [$expectedOld] but found [$current]\"\n    }\n    # Assigning to Cell.Range.Text replaces only the text content while\n    # keeping the existing run formatting (font, size) and paragraph\n    # formatting (justification) on the cell intact.\n    $cell.Range.Text = $newText\n}\n"}
